$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 678.8125
$ws.Range("I41").Value = 898.5
$ws.Range("K41").Value = 898.5
$ws.Range("M41").Value = -458.5

$ws.Range("H100").Value = 3835.8667
$ws.Range("I100").Value = 3854.2
$ws.Range("J100").Value = 3799.2
$ws.Range("K100").Value = 3854.2
$ws.Range("L100").Value = 3799.2
$ws.Range("M100").Value = -3313.2
$ws.Range("N100").Value = -4881.2

$ws.Range("H132").Value = 1506.0735
$ws.Range("I132").Value = 1549.8889
$ws.Range("K132").Value = 4649.6667
$ws.Range("M132").Value = -2119.6667

$ws.Range("H137").Value = 144373.81
$ws.Range("I137").Value = 1875.4
$ws.Range("K137").Value = 5626.200000000001
$ws.Range("M137").Value = -3076.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 347.84616
$ws.Range("I5").Value = 124.77778
$ws.Range("K5").Value = 124.77778
$ws.Range("M5").Value = -12.77778000000001

$ws.Range("H28").Value = 26862.25
$ws.Range("J28").Value = 26499.5
$ws.Range("L28").Value = 26499.5
$ws.Range("N28").Value = -26883.5

$ws.Range("H32").Value = 4982.247
$ws.Range("I32").Value = 3142.9102
$ws.Range("K32").Value = 3142.9102
$ws.Range("M32").Value = -2855.9102

$ws.Range("H74").Value = 493608.38
$ws.Range("I74").Value = 780293.0600000001
$ws.Range("K74").Value = 780293.0600000001
$ws.Range("M74").Value = -779419.0600000001

$ws.Range("H77").Value = 493608.38
$ws.Range("I77").Value = 780293.0600000001
$ws.Range("K77").Value = 3901465.3
$ws.Range("M77").Value = -3897097.3

$ws.Range("H97").Value = 775.25
$ws.Range("I97").Value = 452.73334
$ws.Range("K97").Value = 452.73334
$ws.Range("M97").Value = 43.26666

$ws.Range("H99").Value = 26862.25
$ws.Range("J99").Value = 26499.5
$ws.Range("L99").Value = 26499.5
$ws.Range("N99").Value = -32489.5

$ws.Range("H102").Value = 2764
$ws.Range("I102").Value = 2452
$ws.Range("K102").Value = 2452
$ws.Range("M102").Value = -830

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 347.84616
$ws.Range("I4").Value = 124.77778
$ws.Range("K4").Value = 124.77778
$ws.Range("M4").Value = -9.777780000000007

$ws.Range("H99").Value = 1498.2941
$ws.Range("I99").Value = 1822.7
$ws.Range("J99").Value = 1034.8572
$ws.Range("K99").Value = 1822.7
$ws.Range("L99").Value = 1034.8572
$ws.Range("M99").Value = -324.7
$ws.Range("N99").Value = -4030.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2583.7144
$ws.Range("I31").Value = 2207.0454
$ws.Range("J31").Value = 2785.8293
$ws.Range("K31").Value = 2207.0454
$ws.Range("L31").Value = 2785.8293
$ws.Range("M31").Value = -1912.0454
$ws.Range("N31").Value = -3375.8293

$ws.Range("H34").Value = 2583.7144
$ws.Range("I34").Value = 2207.0454
$ws.Range("J34").Value = 2785.8293
$ws.Range("K34").Value = 2207.0454
$ws.Range("L34").Value = 2785.8293
$ws.Range("M34").Value = -2005.0454
$ws.Range("N34").Value = -3189.8293

$ws.Range("H62").Value = 119333.11
$ws.Range("I62").Value = 202800
$ws.Range("J62").Value = 14999.5
$ws.Range("K62").Value = 202800
$ws.Range("L62").Value = 14999.5
$ws.Range("M62").Value = -202176
$ws.Range("N62").Value = -16247.5

$ws.Range("H65").Value = 119333.11
$ws.Range("I65").Value = 202800
$ws.Range("J65").Value = 14999.5
$ws.Range("K65").Value = 1014000
$ws.Range("L65").Value = 74997.5
$ws.Range("M65").Value = -1010880
$ws.Range("N65").Value = -81237.5

$ws.Range("H69").Value = 33697
$ws.Range("I69").Value = 17394
$ws.Range("K69").Value = 17394
$ws.Range("M69").Value = -16645

$ws.Range("H72").Value = 33697
$ws.Range("I72").Value = 17394
$ws.Range("K72").Value = 52182
$ws.Range("M72").Value = -48438

$ws.Range("H132").Value = 4159.2
$ws.Range("J132").Value = 23840.25
$ws.Range("L132").Value = 71520.75
$ws.Range("N132").Value = -76580.75

$ws.Range("H134").Value = 2799.7673
$ws.Range("I134").Value = 2663.5676
$ws.Range("J134").Value = 3639.6667
$ws.Range("K134").Value = 7990.702799999999
$ws.Range("L134").Value = 10919.0001
$ws.Range("M134").Value = -5455.702799999999
$ws.Range("N134").Value = -15989.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 804.375
$ws.Range("I5").Value = 731.375
$ws.Range("J5").Value = 950.375
$ws.Range("K5").Value = 2194.125
$ws.Range("L5").Value = 2851.125
$ws.Range("M5").Value = -2082.125
$ws.Range("N5").Value = -3075.125

$ws.Range("H34").Value = 592.1429000000001
$ws.Range("I34").Value = 542.6
$ws.Range("J34").Value = 716
$ws.Range("K34").Value = 1627.8
$ws.Range("L34").Value = 2148
$ws.Range("M34").Value = -1543.8
$ws.Range("N34").Value = -2316

$ws.Range("H39").Value = 3866.6667
$ws.Range("J39").Value = 15000
$ws.Range("L39").Value = 45000
$ws.Range("N39").Value = -45588

$ws.Range("H55").Value = 3664.7144
$ws.Range("I55").Value = 1775.5
$ws.Range("K55").Value = 5326.5
$ws.Range("M55").Value = -5149.5

$ws.Range("H107").Value = 1504.6842
$ws.Range("I107").Value = 351.5
$ws.Range("J107").Value = 2036.9231
$ws.Range("K107").Value = 1054.5
$ws.Range("L107").Value = 6110.7693
$ws.Range("M107").Value = 865.5
$ws.Range("N107").Value = -9950.7693

$ws.Range("H135").Value = 804.375
$ws.Range("I135").Value = 731.375
$ws.Range("J135").Value = 950.375
$ws.Range("K135").Value = 6582.375
$ws.Range("L135").Value = 8553.375
$ws.Range("M135").Value = -4047.375
$ws.Range("N135").Value = -13623.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 2803.75
$ws.Range("I23").Value = 2067
$ws.Range("J23").Value = 5014
$ws.Range("K23").Value = 2067
$ws.Range("L23").Value = 5014
$ws.Range("M23").Value = -1844
$ws.Range("N23").Value = -5460

$ws.Range("H70").Value = 9332.333000000001
$ws.Range("I70").Value = 7999
$ws.Range("J70").Value = 10665.667
$ws.Range("K70").Value = 7999
$ws.Range("L70").Value = 10665.667
$ws.Range("M70").Value = -7729
$ws.Range("N70").Value = -11205.667

$ws.Range("H73").Value = 9332.333000000001
$ws.Range("I73").Value = 7999
$ws.Range("J73").Value = 10665.667
$ws.Range("K73").Value = 7999
$ws.Range("L73").Value = 10665.667
$ws.Range("M73").Value = -7063
$ws.Range("N73").Value = -12537.667

$ws.Range("H97").Value = 19250062
$ws.Range("I97").Value = 25024222
$ws.Range("K97").Value = 25024222
$ws.Range("M97").Value = -25023726

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").ClearContents()
$ws.Range("N123").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H60").Value = 39657.8
$ws.Range("J60").Value = 39657.8
$ws.Range("L60").Value = 39657.8
$ws.Range("N60").Value = -40675.8

$ws.Range("H93").Value = 990.8095
$ws.Range("I93").Value = 875.5333000000001
$ws.Range("K93").Value = 875.5333000000001
$ws.Range("M93").Value = 372.4666999999999

$ws.Range("H132").Value = 2535.6538
$ws.Range("I132").Value = 2162.537
$ws.Range("J132").Value = 3375.1667
$ws.Range("K132").Value = 6487.610999999999
$ws.Range("L132").Value = 10125.5001
$ws.Range("M132").Value = -3957.610999999999
$ws.Range("N132").Value = -15185.5001

$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("N134").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6108.2
$ws.Range("I62").Value = 3778.4
$ws.Range("J62").Value = 8438
$ws.Range("K62").Value = 3778.4
$ws.Range("L62").Value = 8438
$ws.Range("M62").Value = -3154.4
$ws.Range("N62").Value = -9686

$ws.Range("H65").Value = 6108.2
$ws.Range("I65").Value = 3778.4
$ws.Range("J65").Value = 8438
$ws.Range("K65").Value = 18892
$ws.Range("L65").Value = 42190
$ws.Range("M65").Value = -15772
$ws.Range("N65").Value = -48430

$ws.Range("H126").Value = 1902.3
$ws.Range("J126").Value = 1458
$ws.Range("L126").Value = 4374
$ws.Range("N126").Value = -9314

$ws.Range("H136").Value = 4589.1055
$ws.Range("I136").Value = 4959.205
$ws.Range("J136").Value = 3787.2222
$ws.Range("K136").Value = 14877.615
$ws.Range("L136").Value = 11361.6666
$ws.Range("M136").Value = -12327.615
$ws.Range("N136").Value = -16461.6666
